# "Agregando ejemplo de estructura Fink"
# Reconfigure the restraint table on the "restric" sheet: the DOF restrained
# at node 1 changes from X (=1) to Y (=2), and a new restraint is added for
# node 4 direction X (previously node 1 direction Y with a -30 prescribed
# displacement is replaced by node 4 direction X with a 0 displacement).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("restric")

$ws.Range("B2").Value = 2
$ws.Range("A3").Value = 4
$ws.Range("B3").Value = 1
$ws.Range("D3").Value = 0

# Make "restric" the active/selected sheet (it was "prop_mat" before),
# with D4 as the active cell.
$ws.Activate()
$ws.Range("D4").Select() | Out-Null
